$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Introduce the new shared strings in the same order the source workbook
# shows them (Column C values for rows 3-5 first, then Column B values),
# so the sharedStrings table ends up in the matching order.
$ws.Range("C3").Value = "Softway21"
$ws.Range("C4").Value = "Softway26"
$ws.Range("C5").Value = "Softway25"

$ws.Range("B3").Value = "Developer"
$ws.Range("B4").Value = "Homeowner"
$ws.Range("B5").Value = "HVAC Dealer"

# B3 keeps the plain/default cell format (no explicit style index),
# matching the source's first-entry formatting quirk.
$ws.Range("B3").Style = "Normal"

# Row 3 - remaining columns
$ws.Range("D3").Value = "Test"
$ws.Range("E3").Value = "Test"
$ws.Range("F3").Value = "Test"
$ws.Range("G3").Value = "Test"
$ws.Range("H3").Value = "Test"
$ws.Range("I3").Value = "Houston"
$ws.Range("J3").Value = "Texas"
$ws.Range("K3").Value = 78479
$ws.Range("L3").Value = 1234567895
$ws.Range("M3").Value = 1234567890
$ws.Range("N3").Value = "rais@softway.com"
$ws.Range("O3").Value = "Test"
$ws.Range("P3").Value = "Claim Status"
$ws.Range("Q3").Value = "Test"
$ws.Range("R3").Value = "This is a test comment..."

# Row 4 - remaining columns
$ws.Range("D4").Value = "Test"
$ws.Range("E4").Value = "Test"
$ws.Range("F4").Value = "Test"
$ws.Range("G4").Value = "Test"
$ws.Range("H4").Value = "Test"
$ws.Range("I4").Value = "Houston"
$ws.Range("J4").Value = "Texas"
$ws.Range("K4").Value = 78479
$ws.Range("L4").Value = 1234567895
$ws.Range("M4").Value = 1234567890
$ws.Range("N4").Value = "rais@softway.com"
$ws.Range("O4").Value = "Test"
$ws.Range("P4").Value = "Claim Status"
$ws.Range("Q4").Value = "Test"
$ws.Range("R4").Value = "This is a test comment..."

# Row 5 - remaining columns
$ws.Range("D5").Value = "Test"
$ws.Range("E5").Value = "Test"
$ws.Range("F5").Value = "Test"
$ws.Range("G5").Value = "Test"
$ws.Range("H5").Value = "Test"
$ws.Range("I5").Value = "Houston"
$ws.Range("J5").Value = "Texas"
$ws.Range("K5").Value = 78479
$ws.Range("L5").Value = 1234567895
$ws.Range("M5").Value = 1234567890
$ws.Range("N5").Value = "rais@softway.com"
$ws.Range("O5").Value = "Test"
$ws.Range("P5").Value = "Claim Status"
$ws.Range("Q5").Value = "Test"
$ws.Range("R5").Value = "This is a test comment..."

# New rows share the same 30pt row height as the existing data row.
$ws.Rows.Item(3).RowHeight = 30
$ws.Rows.Item(4).RowHeight = 30
$ws.Rows.Item(5).RowHeight = 30

$ws.Range("B5").Select()
